# Holly added "S.GISH" as a harvester in bioSamples, so the "harvester"
# column (B) for every data row now shows that name instead of the old
# (incorrect) value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2:B24").Value = "S.GISH"

# Side effects of the edit in the live Excel session: the "harvester"
# column was widened slightly and left selected, and the header row's
# height was normalized to match the rest of the sheet.
$ws.Columns.Item(2).ColumnWidth = 8
$ws.Rows.Item(1).RowHeight = 13.8
$ws.Range("B:B").Select()
